# Applies the "#5: property building done" edit:
#  - sheet1 (土地/land): insert a new land parcel row before the existing one
#  - sheet2 (建物/building): fix up header row + expand each building row to the full schema
#  - sheet3 (汽車/car): add a proper data row (index 32) duplicating the existing record
#  - sheet4 (存款/deposit): insert a new bank-deposit row at the top
#  - sheet5 (股票/stock): insert a new stock holding row (台化) at the top
#  - sheet6 (保險/insurance): add a proper data row (index 102) duplicating the existing record
#  - sheet7 (債務/liabilities): add a proper data row (index 114) duplicating the existing header record

$wb = $excel.ActiveWorkbook

function Set-IndexCellStyle($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------------
# Sheet1: 土地 (land)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("土地")
$ws1.Rows(2).Insert()

Set-IndexCellStyle $ws1.Cells.Item(2,1)
$ws1.Cells.Item(2,1).Value = 14
$ws1.Cells.Item(2,2).Value = "臺北市大安區瑞安段二小段08340000地號"
$ws1.Cells.Item(2,3).Value = 2623
$ws1.Cells.Item(2,4).Value = "10000分之202"
$ws1.Cells.Item(2,5).Value = "賴士葆"
$ws1.Cells.Item(2,6).Value = "83年9月"
$ws1.Cells.Item(2,7).Value = "買賣"
$ws1.Cells.Item(2,8).Value = "(超過五年）"
$ws1.Cells.Item(2,9).Value = "land"
$ws1.Cells.Item(2,10).Value = "normal"
$ws1.Cells.Item(2,11).Value = "2012-04-19"
$ws1.Cells.Item(2,12).Value = "賴士葆"
$ws1.Cells.Item(2,13).Value = 866
$ws1.Cells.Item(2,14).Value = "tmp9edb1"
$ws1.Cells.Item(2,15).Value = 14
$ws1.Cells.Item(2,16).Value = 0.0202
$ws1.Cells.Item(2,17).Value = 52.9846

# ---------------------------------------------------------------------------
# Sheet2: 建物 (building) - was an 8-column sheet where row 1 held live data
# instead of headers; rebuild it fully with the common 17-column schema.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("建物")
$ws2.Cells.Clear()

$headers = @("name","area","share_portion","owner","register_date","register_reason","acquire_value","property_category","category","date","legislator_name","legislator_id","source_file","index","portion","total")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $c = $ws2.Cells.Item(1, $i + 2)
    $c.Value = $headers[$i]
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
    $c.Borders.LineStyle = 1
}

function Fill-BuildingRow($ws, $row, $idx, $name, $area, $share, $owner, $regdate, $regreason, $acqval) {
    Set-IndexCellStyle $ws.Cells.Item($row,1)
    $ws.Cells.Item($row,1).Value = $idx
    $ws.Cells.Item($row,2).Value = $name
    $ws.Cells.Item($row,3).Value = $area
    $ws.Cells.Item($row,4).Value = $share
    $ws.Cells.Item($row,5).Value = $owner
    $ws.Cells.Item($row,6).Value = $regdate
    $ws.Cells.Item($row,7).Value = $regreason
    $ws.Cells.Item($row,8).Value = $acqval
    $ws.Cells.Item($row,9).Value = "land"
    $ws.Cells.Item($row,10).Value = "normal"
    $ws.Cells.Item($row,11).Value = "2012-04-19"
    $ws.Cells.Item($row,12).Value = "賴士葆"
    $ws.Cells.Item($row,13).Value = 866
    $ws.Cells.Item($row,14).Value = "tmp9edb1"
    $ws.Cells.Item($row,15).Value = $idx
    $ws.Cells.Item($row,16).Value = $share
    $ws.Cells.Item($row,17).Value = $area
}

Fill-BuildingRow $ws2 2 20 "臺北市大安區瑞安段二小段02940000建號" 156.6 1 "賴士葆" "83年9月" "買賣" "(超過五年）"
Fill-BuildingRow $ws2 3 21 "臺北市大安區瑞安段二小段02983000建號" 4370.2 0.0178571428571429 "賴士葆" "83年9月" "買賣" "(超過五年地下室停車位）"
Fill-BuildingRow $ws2 4 22 "新北市新店區華城二段00053000建號" 238.91 1 "林良娥" "92年7月" "買賣" "(超堝石年)"

# fix up the two rows whose portion/share text differs from the numeric share ratio
$ws2.Cells.Item(3,4).Value = "56分之1"
$ws2.Cells.Item(3,17).Value = 78.0392857142857

$ws2.Cells.Item(2,16).Value = 1
$ws2.Cells.Item(2,17).Value = 156.6
$ws2.Cells.Item(4,4).Value = "全部"
$ws2.Cells.Item(4,16).Value = 1
$ws2.Cells.Item(4,17).Value = 238.91

# ---------------------------------------------------------------------------
# Sheet3: 汽車 (car)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("汽車")
$ws3.Rows(2).Insert()

Set-IndexCellStyle $ws3.Cells.Item(2,1)
$ws3.Cells.Item(2,1).Value = 32
$ws3.Cells.Item(2,2).Value = "中華休旅車"
$ws3.Cells.Item(2,3).Value = 1997
$ws3.Cells.Item(2,4).Value = "林良娥"
$ws3.Cells.Item(2,5).Value = "94年05月26日"
$ws3.Cells.Item(2,6).Value = "買賣"
$ws3.Cells.Item(2,7).Value = "(超過五年）"

# ---------------------------------------------------------------------------
# Sheet4: 存款 (deposit)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("存款")
$ws4.Rows(2).Insert()

Set-IndexCellStyle $ws4.Cells.Item(2,1)
$ws4.Cells.Item(2,1).Value = 47
$ws4.Cells.Item(2,2).Value = "臺灣土地銀行文山分行"
$ws4.Cells.Item(2,3).Value = "活期儲蓄存款"
$ws4.Cells.Item(2,4).Value = "新臺幣"
$ws4.Cells.Item(2,5).Value = "林良娥"
$ws4.Cells.Item(2,6).Value = 444132

# ---------------------------------------------------------------------------
# Sheet5: 股票 (stock)
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("股票")
$ws5.Rows(2).Insert()

Set-IndexCellStyle $ws5.Cells.Item(2,1)
$ws5.Cells.Item(2,1).Value = 67
$ws5.Cells.Item(2,2).Value = "台化"
$ws5.Cells.Item(2,3).Value = "林良娥"
$ws5.Cells.Item(2,4).Value = 6004
$ws5.Cells.Item(2,5).Value = 10
$ws5.Cells.Item(2,6).Value = "新臺幣"
$ws5.Cells.Item(2,7).Value = 60040
$ws5.Cells.Item(2,8).Value = "stock"
$ws5.Cells.Item(2,9).Value = "normal"
$ws5.Cells.Item(2,10).Value = "2012-04-19"
$ws5.Cells.Item(2,11).Value = "賴士葆"
$ws5.Cells.Item(2,12).Value = 866
$ws5.Cells.Item(2,13).Value = "tmp9edb1"
$ws5.Cells.Item(2,14).Value = 67

# ---------------------------------------------------------------------------
# Sheet6: 保險 (insurance)
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("保險")
$ws6.Rows(2).Insert()

Set-IndexCellStyle $ws6.Cells.Item(2,1)
$ws6.Cells.Item(2,1).Value = 102
$ws6.Cells.Item(2,2).Value = "富邦人壽"
$ws6.Cells.Item(2,3).Value = "安泰人壽靈活理財變額保險甲型"
$ws6.Cells.Item(2,4).Value = "賴士葆"

# ---------------------------------------------------------------------------
# Sheet7: 債務 (liabilities)
# ---------------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("債務")
$ws7.Rows(2).Insert()

Set-IndexCellStyle $ws7.Cells.Item(2,1)
$ws7.Cells.Item(2,1).Value = 114
$ws7.Cells.Item(2,2).Value = "(十二）事業投"
$ws7.Cells.Item(2,3).Value = "資（總金額"
$ws7.Cells.Item(2,4).Value = "新臺幣"
$ws7.Cells.Item(2,5).Value = "元）"
for ($col = 6; $col -le 14; $col++) {
    $ws7.Cells.Item(2, $col).Value = ""
}
